# Updated symbol list with refreshed prices/volumes and a few re-ranked rows (rows 12-16, 41-43).
# Columns D (Price) and E (Volume %) hold numeric-looking values that are stored as TEXT in the
# sheet, not numbers. Assigning a bare numeric/percent-looking string makes Excel silently convert
# the cell to a real Number, which would not match the original text formatting. To keep it text we
# prefix with a literal apostrophe (Excel's standard "treat as text" marker) and then call
# ClearFormats() to drop the quote-prefix style Excel applies when it sees the apostrophe, which
# restores the cell to its original unstyled look (no NumberFormat/quotePrefix left behind).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: BNB
$ws.Range("D2").Value = "'246.58"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'0.60%"
$ws.Range("E2").ClearFormats()

# Row 3: OKB
$ws.Range("D3").Value = "'26.24"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'5.22%"
$ws.Range("E3").ClearFormats()

# Row 4: HuobiToken
$ws.Range("D4").Value = "'5.086"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'0.77%"
$ws.Range("E4").ClearFormats()

# Row 5: Cronos
$ws.Range("D5").Value = "'0.05606"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'-0.28%"
$ws.Range("E5").ClearFormats()

# Row 6: KuCoinToken
$ws.Range("D6").Value = "'6.480"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'-0.87%"
$ws.Range("E6").ClearFormats()

# Row 7: MXToken
$ws.Range("D7").Value = "'0.8134"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'0.56%"
$ws.Range("E7").ClearFormats()

# Row 8: FTXToken
$ws.Range("D8").Value = "'0.8472"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'1.00%"
$ws.Range("E8").ClearFormats()

# Row 9: BitrueCoin
$ws.Range("D9").Value = "'0.02862"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'1.39%"
$ws.Range("E9").ClearFormats()

# Row 10: BitMartToken
$ws.Range("D10").Value = "'0.09387"
$ws.Range("D10").ClearFormats()

# Row 11: BitForexToken
$ws.Range("D11").Value = "'0.001521"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'-0.09%"
$ws.Range("E11").ClearFormats()

# Row 12: TigerCash -> One
$ws.Range("B12").Value = "One"
$ws.Range("C12").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D12").Value = "'0.0005961"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'0.25%"
$ws.Range("E12").ClearFormats()

# Row 13: LEO -> TigerCash
$ws.Range("B13").Value = "TigerCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D13").Value = "'0.006128"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-1.76%"
$ws.Range("E13").ClearFormats()

# Row 14: GateToken -> LEO
$ws.Range("B14").Value = "LEO"
$ws.Range("C14").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D14").Value = "'3.595"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'2.72%"
$ws.Range("E14").ClearFormats()

# Row 15: BTSEToken -> GateToken
$ws.Range("B15").Value = "GateToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D15").Value = "'3.011"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'0.93%"
$ws.Range("E15").ClearFormats()

# Row 16: One -> BTSEToken
$ws.Range("B16").Value = "BTSEToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D16").Value = "'2.055"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'0.08%"
$ws.Range("E16").ClearFormats()

# Row 17: BitpandaEcosystemToken
$ws.Range("E17").Value = "'0.66%"
$ws.Range("E17").ClearFormats()

# Row 18: WazirX
$ws.Range("D18").Value = "'0.1338"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'0.11%"
$ws.Range("E18").ClearFormats()

# Row 19: MandalaExchangeToken
$ws.Range("D19").Value = "'0.07003"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'0.85%"
$ws.Range("E19").ClearFormats()

# Row 20: LiechtensteinCryptoassetsExchange
$ws.Range("E20").Value = "'-2.51%"
$ws.Range("E20").ClearFormats()

# Row 22: MCDex
$ws.Range("D22").Value = "'3.747"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'0.27%"
$ws.Range("E22").ClearFormats()

# Row 23: CoinExToken
$ws.Range("D23").Value = "'0.04649"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'-0.84%"
$ws.Range("E23").ClearFormats()

# Row 24: ZBToken
$ws.Range("E24").Value = "'-1.36%"
$ws.Range("E24").ClearFormats()

# Row 25: BitKan
$ws.Range("D25").Value = "'0.001245"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'0.26%"
$ws.Range("E25").ClearFormats()

# Row 26: HotbitToken
$ws.Range("D26").Value = "'0.004585"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'1.28%"
$ws.Range("E26").ClearFormats()

# Row 27: NitroEx
$ws.Range("D27").Value = "'0.00009601"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'-0.93%"
$ws.Range("E27").ClearFormats()

# Row 28: UpBots
$ws.Range("E28").Value = "'168.03%"
$ws.Range("E28").ClearFormats()

# Row 40: IDEX
$ws.Range("D40").Value = "'0.03669"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'1.08%"
$ws.Range("E40").ClearFormats()

# Row 41: KickToken -> BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1368"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'30.12%"
$ws.Range("E41").ClearFormats()

# Row 42: BKEXToken -> CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002660"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'-2.23%"
$ws.Range("E42").ClearFormats()

# Row 43: CEJI -> KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003397"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'-45.42%"
$ws.Range("E43").ClearFormats()

# Row 44: LocalTraders
$ws.Range("D44").Value = "'0.008773"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'5.30%"
$ws.Range("E44").ClearFormats()

# Row 45: CoinLion
$ws.Range("D45").Value = "'0.00005296"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'0.36%"
$ws.Range("E45").ClearFormats()

# Row 47: CoinbaseStockToken
$ws.Range("D47").Value = "'0.1100"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'-42.06%"
$ws.Range("E47").ClearFormats()

# Row 48: BOLO
$ws.Range("D48").Value = "'0.002663"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'30.05%"
$ws.Range("E48").ClearFormats()

# Row 49: CryptobidCoin
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").ClearFormats()

# Row 50: SpecialPowerGold
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").ClearFormats()
